$wb = $excel.ActiveWorkbook

# Update "想去人数" (interest count) values in both the "展览" and
# "全部类型" sheets, which carry duplicate rows for the same events.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F5").Value = 847
    $ws.Range("F8").Value = 8450
    $ws.Range("F9").Value = 75
    $ws.Range("F12").Value = 111
    $ws.Range("F19").Value = 721
}
